$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$currencyFmt = '_-"$"\ * #,##0_-;\-"$"\ * #,##0_-;_-"$"\ * "-"_-;_-@_-'

# ---------------------------------------------------------------------------
# Row 25: turn the blank trailing row below table "Tabla134" into a Total row
# ---------------------------------------------------------------------------
$ws.Range("A25").Value = "Total"
$ws.Range("C25").NumberFormat = $currencyFmt
$ws.Range("D25").NumberFormat = $currencyFmt
$ws.Range("D25").Formula = "=SUM(D2:D24)"
$ws.ListObjects("Tabla134").Resize($ws.Range("A1:D25"))

# ---------------------------------------------------------------------------
# Row 26: blank spacer row (keep currency formatting consistent with above)
# ---------------------------------------------------------------------------
$ws.Range("C26").NumberFormat = $currencyFmt
$ws.Range("D26").NumberFormat = $currencyFmt

# ---------------------------------------------------------------------------
# Row 27: section title
# ---------------------------------------------------------------------------
$ws.Range("A27").Value = "Componentes adicionales"

# ---------------------------------------------------------------------------
# Row 28: header row for the new table
# ---------------------------------------------------------------------------
$ws.Range("A28").Value = "Producto"
$ws.Range("B28").Value = "Cantidad"
$ws.Range("C28").Value = "Precio"
$ws.Range("D28").Value = "Total"

# ---------------------------------------------------------------------------
# Rows 29-35: new component rows
# ---------------------------------------------------------------------------
$items = @(
    @("PIC16F877A", 1, 4583),
    @("MAX232", 1, 1000),
    @("16x2 LCD", 1, 7668),
    @("Sensor de Luz", 1, 3934),
    @("Sensor de Sonidos", 1, 4900),
    @("Pulsadores", 4, 100),
    @("Cable USB a Serial", 1, 4900)
)

$row = 29
foreach ($item in $items) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).NumberFormat = $currencyFmt
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).NumberFormat = $currencyFmt
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# Create the new table "Tabla1" over the header + data rows.
# NOTE: this engine's COM bridge re-sorts the ListObjects collection by name
# whenever a ListObject.Name is changed, so any variable handle captured
# before the rename can silently refer to the WRONG table afterwards.
# Work around it by always re-resolving tables by their (final) name string.
# ---------------------------------------------------------------------------
$null = $ws.ListObjects.Add(1, $ws.Range("A28:D35"), 0, 1)
$ws.ListObjects.Item(2).Name = "Tabla1"
$ws.ListObjects("Tabla1").TableStyle = "TableStyleMedium6"

$row = 29
foreach ($item in $items) {
    $ws.Cells.Item($row, 4).Formula = "=Tabla1[[#This Row],[Cantidad]]*Tabla1[[#This Row],[Precio]]"
    $row = $row + 1
}

$ws.ListObjects("Tabla1").ShowTotals = $true
$ws.Range("D36").NumberFormat = $currencyFmt
$ws.Range("D36").Formula = "=SUM(Tabla1[Total])"

# ---------------------------------------------------------------------------
# Row 39: grand total combining both tables
# ---------------------------------------------------------------------------
$ws.Range("A39").Value = "Total"
$ws.Range("D39").NumberFormat = $currencyFmt
$ws.Range("D39").Formula = "=D25+Tabla1[[#Totals],[Total]]"

# ---------------------------------------------------------------------------
# Column D width + view
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 12

$ws.Range("E42").Select()
$excel.ActiveWindow.ScrollRow = 21
